# Insert a new weekly data row at row 70. This shifts all the existing
# rows 70..153 down to 71..154 (so the sheet's used range becomes
# A1:R154) and then we populate the freshly inserted row 70 with the
# new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 70..153 down by one row.
$ws.Rows.Item(70).Insert()

# Fill in the new row 70 with the new data point.
$ws.Cells.Item(70, 1).Value  = 6
$ws.Cells.Item(70, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(70, 3).Value  = "Metropolitana"
$ws.Cells.Item(70, 4).Value  = 44571
$ws.Cells.Item(70, 5).Value  = 13
$ws.Cells.Item(70, 6).Value  = 100112001
$ws.Cells.Item(70, 7).Value  = "Berenjena"
$ws.Cells.Item(70, 8).Value  = "Sin especificar"
$ws.Cells.Item(70, 9).Value  = "Primera"
$ws.Cells.Item(70, 10).Value = 370
$ws.Cells.Item(70, 11).Value = 6000
$ws.Cells.Item(70, 12).Value = 7000
$ws.Cells.Item(70, 13).Value = 6541
$ws.Cells.Item(70, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 131
$ws.Cells.Item(70, 17).Value = 50
$ws.Cells.Item(70, 18).Value = "Hortaliza"
